{"js": "// Merge the three separate runs \"<id>\", \"p123r_1\", \"</id>\" (paragraph\n// paraId 00000005) into a single run containing \"<id>p123r_1</id>\".\n// Word's own \"type over a selection\" behaviour keeps the formatting of\n// the first character of the replaced range, so searching for the\n// full combined text and replacing it in place reproduces exactly\n// that merge (first run's rPr survives, trailing empty run untouched).\n\nconst searchResults = context.document.body.search(\"<id>p123r_1</id>\", {\n  matchCase: true,\n  matchWildcards: false,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find text \"<id>p123r_1</id>\" to merge runs.');\n}\n\nconst target = searchResults.items[0];\ntarget.insertText(\"<id>p123r_1</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Merge the three separate runs \"<id>\", \"p123r_1\", \"</id>\" (paragraph\n# paraId 00000005) into a single run containing \"<id>p123r_1</id>\".\n# A Find/Replace over the combined text (Word's Find operates on the\n# paragraph's visible text, independent of run boundaries) collapses\n# the three runs into one, keeping the formatting of the first run\n# (\"<id>\"), exactly like typing the replacement text over the selected\n# range in Word would.\n\n$d = $word.ActiveDocument\n\n$searchText = \"<id>p123r_1</id>\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $searchText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $searchText\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n"}
